# HerbariumHours.xlsx — gui pseudocode & openrefine fix
#
# 1. Fill in the rest of row 8 (week of 10/1-8): B8 becomes 4, and C8/D8/E8
#    get their hours too.
# 2. Re-enter the G2:G7 Total formulas as a single shared formula so they
#    save back out as a shared formula group (master at G2, slaves G3:G7).
# 3. Move the active selection to E9 (just under the newly filled E8).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8 data fixes -------------------------------------------------
$ws.Range("B8").Value = 4
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 2
$ws.Range("E8").Value = 2.5

# --- Shared formula for the Total column (G2:G7) -----------------------
$ws.Range("G2:G7").Formula = "=SUM(B2:F2)"

# --- Selection -----------------------------------------------------------
$ws.Range("E9").Select()
